# This script reproduces a weekly "roll forward" edit on the sheet:
# rows 308..399 are shifted down by 2 rows (to 310..401), and two brand
# new data rows are inserted at rows 308 and 309.
#
# Only columns D, H, I, J, K, L, M, P actually vary row to row in this
# block; columns A, B, C, E, F, G, N, O, Q, R hold constant values for
# every row in the block, so shifting the varying columns down by two
# rows (and writing fresh values into the two newly freed rows) produces
# an identical result to inserting two real rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 308
$lastRowBefore = 399
$shift = 2
$lastRowAfter = $lastRowBefore + $shift

# The two rows being added past the previous end of the sheet (400 and
# 401) do not have a pre-existing style, so writing a date value there
# would otherwise pick up a generic date format. Explicitly copy the
# "Fecha" column's number format so the resulting style matches column
# D everywhere else in the block (numFmtId 165 / style index 2).
$ws.Cells.Item(400, 4).NumberFormat = $ws.Cells.Item(308, 4).NumberFormat()
$ws.Cells.Item(401, 4).NumberFormat = $ws.Cells.Item(308, 4).NumberFormat()

# Walk from the bottom up so we never overwrite a source row before it
# has been read. NOTE: reading a cell's value in this engine must be
# done by calling .Value() (method-call syntax); the bare .Value
# property getter is not resolved to the underlying data.
for ($r = $lastRowBefore; $r -ge $firstRow; $r--) {
    $destRow = $r + $shift
    $ws.Cells.Item($destRow, 4).Value = $ws.Cells.Item($r, 4).Value()    # D - Fecha
    $ws.Cells.Item($destRow, 8).Value = $ws.Cells.Item($r, 8).Value()    # H - Variedad
    $ws.Cells.Item($destRow, 9).Value = $ws.Cells.Item($r, 9).Value()    # I - Calidad
    $ws.Cells.Item($destRow, 10).Value = $ws.Cells.Item($r, 10).Value()  # J - Volumen
    $ws.Cells.Item($destRow, 11).Value = $ws.Cells.Item($r, 11).Value()  # K - Precio minimo
    $ws.Cells.Item($destRow, 12).Value = $ws.Cells.Item($r, 12).Value()  # L - Precio maximo
    $ws.Cells.Item($destRow, 13).Value = $ws.Cells.Item($r, 13).Value()  # M - Precio promedio ponderado
    $ws.Cells.Item($destRow, 16).Value = $ws.Cells.Item($r, 16).Value()  # P - Precio $/Kg
}

# New row 308: Crespo record / Primera
$ws.Cells.Item(308, 4).Value = 44841
$ws.Cells.Item(308, 8).Value = "Crespo record"
$ws.Cells.Item(308, 9).Value = "Primera"
$ws.Cells.Item(308, 10).Value = 1000
$ws.Cells.Item(308, 11).Value = 1500
$ws.Cells.Item(308, 12).Value = 1600
$ws.Cells.Item(308, 13).Value = 1550
$ws.Cells.Item(308, 16).Value = 1550

# New row 309: Crespo record / Segunda
$ws.Cells.Item(309, 4).Value = 44841
$ws.Cells.Item(309, 8).Value = "Crespo record"
$ws.Cells.Item(309, 9).Value = "Segunda"
$ws.Cells.Item(309, 10).Value = 500
$ws.Cells.Item(309, 11).Value = 1100
$ws.Cells.Item(309, 12).Value = 1100
$ws.Cells.Item(309, 13).Value = 1100
$ws.Cells.Item(309, 16).Value = 1100

# The remaining static columns (A, B, C, E, F, G, N, O, Q, R) are
# identical for every row in this block, so they are already correct
# for rows 310..401 (copied from rows 308..399) and for the two new
# rows 308..309 (copied from the existing block values). Fill them in
# explicitly for the two new rows and the two newly extended rows to be
# safe/explicit.
$staticCols = @(1, 2, 3, 5, 6, 7, 14, 15, 17, 18)  # A,B,C,E,F,G,N,O,Q,R
foreach ($col in $staticCols) {
    $val = $ws.Cells.Item(307, $col).Value()
    $ws.Cells.Item(308, $col).Value = $val
    $ws.Cells.Item(309, $col).Value = $val
    $ws.Cells.Item(400, $col).Value = $val
    $ws.Cells.Item(401, $col).Value = $val
}
